# Updates Sema3d-Nrp1 LR-pairs data (columns E:T) with recomputed TPM-based values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.439567666666667
$ws.Range("H2").Value = 7.318703
$ws.Range("I2").Value = 0.8306928434252403
$ws.Range("J2").Value = 0.8306928434252402
$ws.Range("M2").Value = 133.7780026666667
$ws.Range("N2").Value = 401.334008
$ws.Range("O2").Value = 0.50863533211804
$ws.Range("P2").Value = 0.5086353321180399
$ws.Range("Q2").Value = 326.3604898168471
$ws.Range("R2").Value = 2937.244408351624
$ws.Range("S2").Value = 0.4225197303036761
$ws.Range("T2").Value = 0.422519730303676
# Row 3
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.439567666666667
$ws.Range("H3").Value = 7.318703
$ws.Range("I3").Value = 0.8306928434252403
$ws.Range("J3").Value = 0.8306928434252402
$ws.Range("O3").Value = 0.1993888292903622
$ws.Range("P3").Value = 0.1993888292903622
$ws.Range("Q3").Value = 127.935736827871
$ws.Range("R3").Value = 1151.421631450839
$ws.Range("S3").Value = 0.1656308735504408
$ws.Range("T3").Value = 0.1656308735504408
# Row 4
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.439567666666667
$ws.Range("H4").Value = 7.318703
$ws.Range("I4").Value = 0.8306928434252403
$ws.Range("J4").Value = 0.8306928434252402
$ws.Range("M4").Value = 21.197691
$ws.Range("N4").Value = 63.593073
$ws.Range("O4").Value = 0.08059542216956049
$ws.Range("P4").Value = 0.08059542216956046
$ws.Range("Q4").Value = 51.71320157159101
$ws.Range("R4").Value = 465.418814144319
$ws.Range("S4").Value = 0.06695004040908985
$ws.Range("T4").Value = 0.06695004040908982
# Row 5
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.439567666666667
$ws.Range("H5").Value = 7.318703
$ws.Range("I5").Value = 0.8306928434252403
$ws.Range("J5").Value = 0.8306928434252402
$ws.Range("M5").Value = 55.59592133333333
$ws.Range("N5").Value = 166.787764
$ws.Range("O5").Value = 0.2113804164220374
$ws.Range("P5").Value = 0.2113804164220373
$ws.Range("Q5").Value = 135.6300120833435
$ws.Range("R5").Value = 1220.670108750092
$ws.Range("S5").Value = 0.1755921991620336
$ws.Range("T5").Value = 0.1755921991620335
# Row 6
$ws.Range("G6").Value = 0.47512
$ws.Range("H6").Value = 1.42536
$ws.Range("I6").Value = 0.1617822654238874
$ws.Range("J6").Value = 0.1617822654238873
$ws.Range("M6").Value = 133.7780026666667
$ws.Range("N6").Value = 401.334008
$ws.Range("O6").Value = 0.50863533211804
$ws.Range("P6").Value = 0.5086353321180399
$ws.Range("Q6").Value = 63.56060462698667
$ws.Range("R6").Value = 572.0454416428801
$ws.Range("S6").Value = 0.08228817630468785
$ws.Range("T6").Value = 0.0822881763046878
# Row 7
$ws.Range("G7").Value = 0.47512
$ws.Range("H7").Value = 1.42536
$ws.Range("I7").Value = 0.1617822654238874
$ws.Range("J7").Value = 0.1617822654238873
$ws.Range("O7").Value = 0.1993888292903622
$ws.Range("P7").Value = 0.1993888292903622
$ws.Range("Q7").Value = 24.91622926152
$ws.Range("R7").Value = 224.24606335368
$ws.Range("S7").Value = 0.03225757650281155
$ws.Range("T7").Value = 0.03225757650281153
# Row 8
$ws.Range("G8").Value = 0.47512
$ws.Range("H8").Value = 1.42536
$ws.Range("I8").Value = 0.1617822654238874
$ws.Range("J8").Value = 0.1617822654238873
$ws.Range("M8").Value = 21.197691
$ws.Range("N8").Value = 63.593073
$ws.Range("O8").Value = 0.08059542216956049
$ws.Range("P8").Value = 0.08059542216956046
$ws.Range("Q8").Value = 10.07144694792
$ws.Range("R8").Value = 90.64302253128
$ws.Range("S8").Value = 0.01303890998138609
$ws.Range("T8").Value = 0.01303890998138608
# Row 9
$ws.Range("G9").Value = 0.47512
$ws.Range("H9").Value = 1.42536
$ws.Range("I9").Value = 0.1617822654238874
$ws.Range("J9").Value = 0.1617822654238873
$ws.Range("M9").Value = 55.59592133333333
$ws.Range("N9").Value = 166.787764
$ws.Range("O9").Value = 0.2113804164220374
$ws.Range("P9").Value = 0.2113804164220373
$ws.Range("Q9").Value = 26.41473414389333
$ws.Range("R9").Value = 237.73260729504
$ws.Range("S9").Value = 0.03419760263500188
$ws.Range("T9").Value = 0.03419760263500186
# Row 10
$ws.Range("E10").Value = 1
$ws.Range("F10").Value = 0.3333333333333333
$ws.Range("G10").Value = 0.022099
$ws.Range("H10").Value = 0.066297
$ws.Range("I10").Value = 0.007524891150872382
$ws.Range("J10").Value = 0.00752489115087238
$ws.Range("M10").Value = 133.7780026666667
$ws.Range("N10").Value = 401.334008
$ws.Range("O10").Value = 0.50863533211804
$ws.Range("P10").Value = 0.5086353321180399
$ws.Range("Q10").Value = 2.956360080930667
$ws.Range("R10").Value = 26.607240728376
$ws.Range("S10").Value = 0.003827425509676075
$ws.Range("T10").Value = 0.003827425509676073
# Row 11
$ws.Range("E11").Value = 1
$ws.Range("F11").Value = 0.3333333333333333
$ws.Range("G11").Value = 0.022099
$ws.Range("H11").Value = 0.066297
$ws.Range("I11").Value = 0.007524891150872382
$ws.Range("J11").Value = 0.00752489115087238
$ws.Range("O11").Value = 0.1993888292903622
$ws.Range("P11").Value = 0.1993888292903622
$ws.Range("Q11").Value = 1.158915117129
$ws.Range("R11").Value = 10.430236054161
$ws.Range("S11").Value = 0.001500379237109851
$ws.Range("T11").Value = 0.00150037923710985
# Row 12
$ws.Range("E12").Value = 1
$ws.Range("F12").Value = 0.3333333333333333
$ws.Range("G12").Value = 0.022099
$ws.Range("H12").Value = 0.066297
$ws.Range("I12").Value = 0.007524891150872382
$ws.Range("J12").Value = 0.00752489115087238
$ws.Range("M12").Value = 21.197691
$ws.Range("N12").Value = 63.593073
$ws.Range("O12").Value = 0.08059542216956049
$ws.Range("P12").Value = 0.08059542216956046
$ws.Range("Q12").Value = 0.468447773409
$ws.Range("R12").Value = 4.216029960681
$ws.Range("S12").Value = 0.0006064717790845495
$ws.Range("T12").Value = 0.0006064717790845492
# Row 13
$ws.Range("E13").Value = 1
$ws.Range("F13").Value = 0.3333333333333333
$ws.Range("G13").Value = 0.022099
$ws.Range("H13").Value = 0.066297
$ws.Range("I13").Value = 0.007524891150872382
$ws.Range("J13").Value = 0.00752489115087238
$ws.Range("M13").Value = 55.59592133333333
$ws.Range("N13").Value = 166.787764
$ws.Range("O13").Value = 0.2113804164220374
$ws.Range("P13").Value = 0.2113804164220373
$ws.Range("Q13").Value = 1.228614265545333
$ws.Range("R13").Value = 11.057528389908
$ws.Range("S13").Value = 0.001590614625001908
$ws.Range("T13").Value = 0.001590614625001907
